# Update header label for column B
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Value = "Active Cases"

# Update/add Cluster Name + Active Cases rows (A2:B70)
$ws.Cells.Item(2, 1).Value = "3398 BlueCross Elly Kay Mordialloc"
$ws.Cells.Item(2, 2).Value = 29
$ws.Cells.Item(3, 1).Value = "3564 Waverley Valley Aged Care Glen Waverley"
$ws.Cells.Item(3, 2).Value = 15
$ws.Cells.Item(4, 1).Value = "3601 Baptcare Westhaven community"
$ws.Cells.Item(4, 2).Value = 13
$ws.Cells.Item(5, 1).Value = "3647 Aurrum Aged Care Reservoir"
$ws.Cells.Item(5, 2).Value = 11
$ws.Cells.Item(6, 1).Value = "3653 Fronditha Thalpori St Albans Aged Care"
$ws.Cells.Item(6, 2).Value = 20
$ws.Cells.Item(7, 1).Value = "3975 Aurrum Aged Care Brunswick West"
$ws.Cells.Item(7, 2).Value = 13
$ws.Cells.Item(8, 1).Value = "4257 BlueCross The Gables Camberwell"
$ws.Cells.Item(8, 2).Value = 16
$ws.Cells.Item(9, 1).Value = "4295 Hope Aged Care Sunshine West"
$ws.Cells.Item(9, 2).Value = 16
$ws.Cells.Item(10, 1).Value = "4314 Estia Health Ardeer"
$ws.Cells.Item(10, 2).Value = 17
$ws.Cells.Item(11, 1).Value = "44095 Myrniong Primary School Myrniong"
$ws.Cells.Item(11, 2).Value = 13
$ws.Cells.Item(12, 1).Value = "44304 Brighton Primary School Brighton"
$ws.Cells.Item(12, 2).Value = 14
$ws.Cells.Item(13, 1).Value = "44404 Castlemaine North Primary School Castlemaine"
$ws.Cells.Item(13, 2).Value = 41
$ws.Cells.Item(14, 1).Value = "44490 Armadale Primary School Armadale"
$ws.Cells.Item(14, 2).Value = 27
$ws.Cells.Item(15, 1).Value = "44593 Torquay P-6 College Torquay"
$ws.Cells.Item(15, 2).Value = 34
$ws.Cells.Item(16, 1).Value = "44620 Canterbury Primary School Canterbury"
$ws.Cells.Item(16, 2).Value = 16
$ws.Cells.Item(17, 1).Value = "44623 Brunswick North Primary School Brunswick West"
$ws.Cells.Item(17, 2).Value = 27
$ws.Cells.Item(18, 1).Value = "44745 Briar Hill Primary School Briar Hill"
$ws.Cells.Item(18, 2).Value = 22
$ws.Cells.Item(19, 1).Value = "44799 Eastwood Primary School Ringwood East"
$ws.Cells.Item(19, 2).Value = 37
$ws.Cells.Item(20, 1).Value = "44960 Thomastown West Primary School"
$ws.Cells.Item(20, 2).Value = 14
$ws.Cells.Item(21, 1).Value = "45013 Gladstone Views Primary School"
$ws.Cells.Item(21, 2).Value = 23
$ws.Cells.Item(22, 1).Value = "45147 Maramba Primary School Narre Warren"
$ws.Cells.Item(22, 2).Value = 13
$ws.Cells.Item(23, 1).Value = "45168 Ranfurly Primary School Mildura"
$ws.Cells.Item(23, 2).Value = 21
$ws.Cells.Item(24, 1).Value = "45181 Courtenay Gardens Primary School Cranbourne North"
$ws.Cells.Item(24, 2).Value = 10
$ws.Cells.Item(25, 1).Value = "45226 Glen Waverley South Primary School Glen Waverley"
$ws.Cells.Item(25, 2).Value = 10
$ws.Cells.Item(26, 1).Value = "45257 Roxburgh Rise Primary School Roxburgh Park"
$ws.Cells.Item(26, 2).Value = 12
$ws.Cells.Item(27, 1).Value = "45305 Lockington Consolidated School Lockington"
$ws.Cells.Item(27, 2).Value = 28
$ws.Cells.Item(28, 1).Value = "45719 St Joseph's Primary School Numurkah"
$ws.Cells.Item(28, 2).Value = 11
$ws.Cells.Item(29, 1).Value = "4574 Village Glen Aged Care Residences Mornington"
$ws.Cells.Item(29, 2).Value = 10
$ws.Cells.Item(30, 1).Value = "45757 Saint Joseph's Primary School Warragul"
$ws.Cells.Item(30, 2).Value = 12
$ws.Cells.Item(31, 1).Value = "45764 Our Lady Help of Christian's Primary School Brunswick East"
$ws.Cells.Item(31, 2).Value = 11
$ws.Cells.Item(32, 1).Value = "45858 St Bernard's Primary Coburg"
$ws.Cells.Item(32, 2).Value = 28
$ws.Cells.Item(33, 1).Value = "45861 St Oliver Plunkett Primary School Pascoe Vale"
$ws.Cells.Item(33, 2).Value = 10
$ws.Cells.Item(34, 1).Value = "45958 Ave Maria College Aberfeldie Workplace"
$ws.Cells.Item(34, 2).Value = 20
$ws.Cells.Item(35, 1).Value = "45975 St Thomas More Primary School Hadfield"
$ws.Cells.Item(35, 2).Value = 14
$ws.Cells.Item(36, 1).Value = "46074 St Justin's Catholic Primary School Wheelers Hill"
$ws.Cells.Item(36, 2).Value = 15
$ws.Cells.Item(37, 1).Value = "46078 Corpus Christi Primary School Werribee"
$ws.Cells.Item(37, 2).Value = 34
$ws.Cells.Item(38, 1).Value = "46086 St Kevin's Primary School Hampton Park"
$ws.Cells.Item(38, 2).Value = 11
$ws.Cells.Item(39, 1).Value = "46104 Clairvaux Catholic School Belmont"
$ws.Cells.Item(39, 2).Value = 10
$ws.Cells.Item(40, 1).Value = "46135 Wesley College Junior School St Kilda Road Melbourne"
$ws.Cells.Item(40, 2).Value = 10
$ws.Cells.Item(41, 1).Value = "46208 Mount Scopus Memorial College Gandel Campus Burwood"
$ws.Cells.Item(41, 2).Value = 13
$ws.Cells.Item(42, 1).Value = "46327 Victory Christian College Strathdale"
$ws.Cells.Item(42, 2).Value = 15
$ws.Cells.Item(43, 1).Value = "50279 Dallas Brooks Community Primary School Dallas"
$ws.Cells.Item(43, 2).Value = 12
$ws.Cells.Item(44, 1).Value = "50584 St Mary of the Cross MacKillop Primary School Epping"
$ws.Cells.Item(44, 2).Value = 11
$ws.Cells.Item(45, 1).Value = "51529 Sirius College Primary School Dallas"
$ws.Cells.Item(45, 2).Value = 11
$ws.Cells.Item(46, 1).Value = "52390 Our Lady of the Way Catholic Primary School Wallan"
$ws.Cells.Item(46, 2).Value = 16
$ws.Cells.Item(47, 1).Value = "52694 Pakenham Primary School Pakenham"
$ws.Cells.Item(47, 2).Value = 14
$ws.Cells.Item(48, 1).Value = "Australian Radio Network Richmond"
$ws.Cells.Item(48, 2).Value = 14
$ws.Cells.Item(49, 1).Value = "Ballarat Freedom Protest"
$ws.Cells.Item(49, 2).Value = 10
$ws.Cells.Item(50, 1).Value = "Brandon Park Primary School Wheelers Hill"
$ws.Cells.Item(50, 2).Value = 13
$ws.Cells.Item(51, 1).Value = "Cardinia Waters Retirement Village Pakenham"
$ws.Cells.Item(51, 2).Value = 13
$ws.Cells.Item(52, 1).Value = "Chisholm Road Prison Project Lara"
$ws.Cells.Item(52, 2).Value = 18
$ws.Cells.Item(53, 1).Value = "Confirmed Omicron Sircuit Bar Fitzroy"
$ws.Cells.Item(53, 2).Value = 14
$ws.Cells.Item(54, 1).Value = "Confirmed Omicron Variant The Peel Hotel Collingwood"
$ws.Cells.Item(54, 2).Value = 17
$ws.Cells.Item(55, 1).Value = "Gladstone Views Primary School Gladstone Park"
$ws.Cells.Item(55, 2).Value = 14
$ws.Cells.Item(56, 1).Value = "Goodstart Early Learning Preston"
$ws.Cells.Item(56, 2).Value = 13
$ws.Cells.Item(57, 1).Value = "Greendale Hotel Greendale"
$ws.Cells.Item(57, 2).Value = 10
$ws.Cells.Item(58, 1).Value = "JBS Australia Brooklyn"
$ws.Cells.Item(58, 2).Value = 29
$ws.Cells.Item(59, 1).Value = "Kororoit Creek Primary School Burnside Heights"
$ws.Cells.Item(59, 2).Value = 26
$ws.Cells.Item(60, 1).Value = "Mambourin Enterprises Allara Deer Park"
$ws.Cells.Item(60, 2).Value = 23
$ws.Cells.Item(61, 1).Value = "Oakleigh South Primary School Oakleigh South"
$ws.Cells.Item(61, 2).Value = 15
$ws.Cells.Item(62, 1).Value = "PGL Camp Rumbug Foster North"
$ws.Cells.Item(62, 2).Value = 47
$ws.Cells.Item(63, 1).Value = "Rosebud Primary School Rosebud"
$ws.Cells.Item(63, 2).Value = 20
$ws.Cells.Item(64, 1).Value = "St Christophers Primary School Airport West"
$ws.Cells.Item(64, 2).Value = 10
$ws.Cells.Item(65, 1).Value = "St Mary's Parish Primary School"
$ws.Cells.Item(65, 2).Value = 11
$ws.Cells.Item(66, 1).Value = "St Vincents Hospital Melbourne Emergency Department Fitzroy"
$ws.Cells.Item(66, 2).Value = 11
$ws.Cells.Item(67, 1).Value = "StarTrack Tullamarine"
$ws.Cells.Item(67, 2).Value = 20
$ws.Cells.Item(68, 1).Value = "The George Lounge St Kilda"
$ws.Cells.Item(68, 2).Value = 15
$ws.Cells.Item(69, 1).Value = "Thomastown West Primary School Camp Doxa's Malmsbury"
$ws.Cells.Item(69, 2).Value = 20
$ws.Cells.Item(70, 1).Value = "V & G construction site San Lorenzo Wine & Dining"
$ws.Cells.Item(70, 2).Value = 14
